$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the "Recorded By" column dynamically from the header row (row 1).
$usedCols = $ws.UsedRange.Columns.Count
$col = 0
for ($c = 1; $c -le $usedCols; $c++) {
    $h = $ws.Cells.Item(1, $c).Value()
    if ($null -ne $h -and [string]$h -eq 'Recorded By') {
        $col = $c
        break
    }
}
if ($col -eq 0) { $col = 7 }

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value()
    if ($null -eq $val) { continue }
    $s = [string]$val
    if ($s -eq '') { continue }
    if ($s -notmatch 'System') { continue }

    # Split the "Recorded By" list on commas (values look like "a@b.com, System").
    $parts = $s -split ',\s*'

    # Find the (case-sensitive) "System" token - PowerShell's -eq/-ceq operators in
    # this environment are not reliably case-sensitive, so use String.Equals instead.
    $idx = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i].Equals('System')) {
            $idx = $i
            break
        }
    }
    if ($idx -lt 0) { continue }

    $rest = @()
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($i -ne $idx) { $rest += $parts[$i] }
    }

    # Put "System" right after a leading lowercase "system" token, otherwise put it first.
    if ($rest.Length -gt 0 -and $rest[0].Equals('system')) {
        if ($rest.Length -gt 1) {
            $newParts = @($rest[0], 'System') + $rest[1..($rest.Length - 1)]
        } else {
            $newParts = @($rest[0], 'System')
        }
    } else {
        $newParts = @('System') + $rest
    }

    $newVal = [string]::Join(', ', $newParts)
    if ($newVal -ne $s) {
        $cell.Value = $newVal
    }
}
